$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 449; this shifts rows 449:511 down to 450:512
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with its data
$ws.Range("A449").Value = 4
$ws.Range("B449").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C449").Value = "Los Lagos"
$ws.Range("D449").Value = 44951
$ws.Range("E449").Value = 10
$ws.Range("F449").Value = "Fruta"
$ws.Range("G449").Value = 100102
$ws.Range("H449").Value = "Cítricos"
$ws.Range("I449").Value = 100102006
$ws.Range("J449").Value = "Pomelo"
$ws.Range("K449").Value = "Start Ruby"
$ws.Range("L449").Value = "Primera"
$ws.Range("M449").Value = 200
$ws.Range("N449").Value = 11000
$ws.Range("O449").Value = 12000
$ws.Range("P449").Value = 11500
$ws.Range("Q449").Value = "$/caja 14 kilos empedrada"
$ws.Range("R449").Value = "Región de O'Higgins"
$ws.Range("S449").Value = 821
$ws.Range("T449").Value = 14
